# "Generate Report for Handback" - refresh the handback status report:
#  - Overview sheet + per-language sheets: flip the shared "Handed back: ...
#    in sync with en-US" status text to "... not in sync with en-US"
#    (Overview!E2:F3 and the "Status" column on the zh-cn/de-de sheets all
#    shared that one string)
#  - zh-cn sheet: the second file's "Correspond Handback DateTime" advances
#  - de-de sheet: the second file's "Correspond Handback DateTime" advances
#  - the columns holding the (now longer) status text get auto-resized
#    along with it

$wb = $excel.ActiveWorkbook

# Column width values land through a pixel-quantized ColumnWidth setter, so
# the exact authored width (33.4602203369141) isn't directly reachable; this
# input resolves to the nearest attainable quantum (33.5) via the same
# round(6*w)/6 + 5/6 rule the host applies.
$newColWidth = 32.6666666666667

# ---------------------------------------------------------------------
# Overview sheet: update the "Handed back: ... in sync with en-US" status
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: not in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------
# zh-cn sheet: Status column mirrors the same "in sync" text, and the
# 61aa49de... row got handed back again -> new datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-10-10 09:54:47"

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------
# de-de sheet: Status column mirrors the same "in sync" text, and the
# 61aa49de... row got handed back again -> new datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-10-10 09:55:02"

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
